# Maplink.xlsx edit: replace the old "Port Blair tourism" map-link list with a
# short "art gallery" map-link list (3 places + a repeated category label),
# turn the first link into a real hyperlink, and tidy up the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urlGallery  = "https://www.google.co.in/maps/place/My+Art+Gallary/data=!4m7!3m6!1s0x390ce5e7cce45dd5:0xad1c1016086fd11b!8m2!3d28.5821195!4d77.3266991!16s%2Fg%2F11rst3b1yv!19sChIJ1V3kzOflDDkRG9FvCBYQHK0?authuser=0&hl=en&rclk=1"
$urlJapingka = "https://www.google.co.in/maps/place/Japingka+Aboriginal+Art/data=!4m7!3m6!1s0x2a32a170196dc9c7:0x1d72bcd4480b9178!8m2!3d-32.055214!4d115.7444264!16s%2Fg%2F1v6gbjsr!19sChIJx8ltGXChMioReJELSNS8ch0?authuser=0&hl=en&rclk=1"
$urlRebecca  = "https://www.google.co.in/maps/place/Rebecca+Hossack+Art+Gallery/data=!4m7!3m6!1s0x48761b2c134a2431:0x5a72d39e5c1bed17!8m2!3d51.5220363!4d-0.139547!16s%2Fg%2F1td5wn6m!19sChIJMSRKEywbdkgRF-0bXJ7Tclo?authuser=0&hl=en&rclk=1"
$label       = "Aboriginal art gallery"

# Drop the old 15-row list, keep only the 4 rows we need.
$ws.Range("A5:A15").ClearContents() | Out-Null

# Column A: the three map links (Japingka repeated on rows 2 and 3, like the source data).
$ws.Range("A1").Value = $urlGallery
$ws.Range("A2").Value = $urlJapingka
$ws.Range("A3").Value = $urlJapingka
$ws.Range("A4").Value = $urlRebecca

# Column B: repeated category label for every row.
$ws.Range("B1").Value = $label
$ws.Range("B2").Value = $label
$ws.Range("B3").Value = $label
$ws.Range("B4").Value = $label

# Make the first link a real (clickable) hyperlink - this also wires up the
# "Hyperlink" cell style/theme color + underline automatically.
$ws.Hyperlinks.Add($ws.Range("A1"), $urlGallery) | Out-Null

# Column A was widened (best-fit) to fit the long URLs.
$ws.Columns.Item(1).ColumnWidth = 227.5

# Leave the selection on the newly-added label column, like the saved file.
$ws.Range("B1:B4").Select() | Out-Null

# Strip personal/path info from the workbook (Inspect Document > Remove Personal
# Information) - corresponds to workbookPr/@filterPrivacy and dropping the
# cached absolute file path in the saved XML.
$wb.RemovePersonalInformation = $true

# Restore the (maximized) window size recorded in the saved file.
$win = $wb.Windows.Item(1)
$win.Width = 22260
$win.Height = 12645

Write-Host "Maplink.xlsx updated: $($ws.UsedRange.Address()) now holds the art-gallery links."
